# Generate Report for Handback
# Adds a new handback row (a278e830-8c0f-430d-b262-e3153c0360f6.md) to all
# three worksheets (Overview, zh-cn, de-de), while also refreshing the
# previously-pending row (81ff3386-...) to its completed state
# (871c2a8f-e428-4c32-a3d7-4f079e4772ac.md) with new hashes/timestamps.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# New / updated identifiers used across the three sheets
# ---------------------------------------------------------------------
$file1 = "871c2a8f-e428-4c32-a3d7-4f079e4772ac.md"
$file2 = "a278e830-8c0f-430d-b262-e3153c0360f6.md"

$zhHash1 = "871c2a8f-e428-4c32-a3d7-4f079e4772ac.f50b0bc230aa48e54d32b4d2ccdadeed5a020b4a.zh-cn.xlf"
$zhHash2 = "a278e830-8c0f-430d-b262-e3153c0360f6.8788f3273d79f0ab53795567bf8e8702f5add239.zh-cn.xlf"

$deHash1 = "871c2a8f-e428-4c32-a3d7-4f079e4772ac.f50b0bc230aa48e54d32b4d2ccdadeed5a020b4a.de-de.xlf"
$deHash2 = "a278e830-8c0f-430d-b262-e3153c0360f6.8788f3273d79f0ab53795567bf8e8702f5add239.de-de.xlf"

$hoDate = "2016-08-17 10:58:10"

$zhHoDate1 = "2016-08-17 10:58:00"
$zhHbDate1 = "2016-08-17 10:58:29"
$zhHoDate2 = $zhHoDate1
$zhHbDate2 = $zhHbDate1

$deHbDate1 = "2016-08-17 10:58:36"
$deHbDate2 = $deHbDate1

$statusText = "Handed back: in sync with en-US"
$ext = ".md"
$e2e = "e2e"
$ht = "ht"

# ===========================================================================
# Sheet 1: Overview
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

# refresh existing (row 2) + add new (row 3) hyperlinks cleanly
$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value = $file1
$ws1.Range("B2").Value = "$e2e\$file1"
$ws1.Range("C2").Value = $ext
$ws1.Range("E2").Value = $statusText
$ws1.Range("F2").Value = $statusText
$ws1.Range("G2").Value = $hoDate
$ws1.Range("G2").NumberFormat = $dateFmt

$ws1.Range("A3").Value = $file2
$ws1.Range("B3").Value = "$e2e\$file2"
$ws1.Range("C3").Value = $ext
$ws1.Range("E3").Value = $statusText
$ws1.Range("F3").Value = $statusText
$ws1.Range("G3").Value = $hoDate
$ws1.Range("G3").NumberFormat = $dateFmt

$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5f986634baa7eba632b415794ab4d209de27a33/e2e/$file1", "", "", "$e2e\$file1") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5f986634baa7eba632b415794ab4d209de27a33/e2e/$file2", "", "", "$e2e\$file2") | Out-Null

# ===========================================================================
# Sheet 2: zh-cn
# ===========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value = $file1
$ws2.Range("B2").Value = $ext
$ws2.Range("C2").Value = $statusText
$ws2.Range("D2").Value = $e2e
$ws2.Range("E2").Value = $ht
$ws2.Range("F2").Value = "'False"
$ws2.Range("G2").Value = $zhHash1
$ws2.Range("H2").Value = $zhHoDate1
$ws2.Range("H2").NumberFormat = $dateFmt
$ws2.Range("I2").Value = $file1
$ws2.Range("J2").Value = $zhHash1
$ws2.Range("K2").Value = $zhHbDate1
$ws2.Range("K2").NumberFormat = $dateFmt
$ws2.Range("L2").Value = "'"
$ws2.Range("M2").Value = "'True"
$ws2.Range("N2").Value = "'"
$ws2.Range("O2").Value = "'False"
$ws2.Range("P2").Value = "'"

$ws2.Range("A3").Value = $file2
$ws2.Range("B3").Value = $ext
$ws2.Range("C3").Value = $statusText
$ws2.Range("D3").Value = $e2e
$ws2.Range("E3").Value = $ht
$ws2.Range("F3").Value = "'True"
$ws2.Range("G3").Value = $zhHash2
$ws2.Range("H3").Value = $zhHoDate2
$ws2.Range("H3").NumberFormat = $dateFmt
$ws2.Range("I3").Value = $file2
$ws2.Range("J3").Value = $zhHash2
$ws2.Range("K3").Value = $zhHbDate2
$ws2.Range("K3").NumberFormat = $dateFmt
$ws2.Range("L3").Value = "'"
$ws2.Range("M3").Value = "'True"
$ws2.Range("N3").Value = "'"
$ws2.Range("O3").Value = "'False"
$ws2.Range("P3").Value = "'"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5f986634baa7eba632b415794ab4d209de27a33/e2e/$file1", "", "", $file1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9c8ea3a7e8519ea27fb8ce761243698b66f4b4d9/e2e/$file1", "", "", $file1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5f986634baa7eba632b415794ab4d209de27a33/e2e/$file2", "", "", $file2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9c8ea3a7e8519ea27fb8ce761243698b66f4b4d9/e2e/$file2", "", "", $file2) | Out-Null

# ===========================================================================
# Sheet 3: de-de
# ===========================================================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value = $file1
$ws3.Range("B2").Value = $ext
$ws3.Range("C2").Value = $statusText
$ws3.Range("D2").Value = $e2e
$ws3.Range("E2").Value = $ht
$ws3.Range("F2").Value = "'False"
$ws3.Range("G2").Value = $deHash1
$ws3.Range("H2").Value = $hoDate
$ws3.Range("H2").NumberFormat = $dateFmt
$ws3.Range("I2").Value = $file1
$ws3.Range("J2").Value = $deHash1
$ws3.Range("K2").Value = $deHbDate1
$ws3.Range("K2").NumberFormat = $dateFmt
$ws3.Range("L2").Value = "'"
$ws3.Range("M2").Value = "'True"
$ws3.Range("N2").Value = "'"
$ws3.Range("O2").Value = "'False"
$ws3.Range("P2").Value = "'"

$ws3.Range("A3").Value = $file2
$ws3.Range("B3").Value = $ext
$ws3.Range("C3").Value = $statusText
$ws3.Range("D3").Value = $e2e
$ws3.Range("E3").Value = $ht
$ws3.Range("F3").Value = "'True"
$ws3.Range("G3").Value = $deHash2
$ws3.Range("H3").Value = $hoDate
$ws3.Range("H3").NumberFormat = $dateFmt
$ws3.Range("I3").Value = $file2
$ws3.Range("J3").Value = $deHash2
$ws3.Range("K3").Value = $deHbDate2
$ws3.Range("K3").NumberFormat = $dateFmt
$ws3.Range("L3").Value = "'"
$ws3.Range("M3").Value = "'True"
$ws3.Range("N3").Value = "'"
$ws3.Range("O3").Value = "'False"
$ws3.Range("P3").Value = "'"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5f986634baa7eba632b415794ab4d209de27a33/e2e/$file1", "", "", $file1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d1e4950d771b0f21abbd676d9d7f000e2a0265de/e2e/$file1", "", "", $file1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5f986634baa7eba632b415794ab4d209de27a33/e2e/$file2", "", "", $file2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d1e4950d771b0f21abbd676d9d7f000e2a0265de/e2e/$file2", "", "", $file2) | Out-Null

# ===========================================================================
# Extend the tables (ListObjects) on each sheet so headerRowCount / ref /
# autoFilter cover the newly populated row 3.
# ===========================================================================
foreach ($ws in @($ws1, $ws2, $ws3)) {
    if ($ws.ListObjects.Count -gt 0) {
        $tbl = $ws.ListObjects.Item(1)
        $lastRow = $tbl.Range.Rows.Count
        if ($lastRow -lt 3) {
            $tbl.Resize($ws.Range($tbl.Range.Cells.Item(1,1), $ws.Cells.Item(3, $tbl.Range.Columns.Count)))
        }
    }
}

Write-Host "Handback report updated."
